$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark tasks 1.3, 1.4, 1.5 (rows 5-7) as "done" in the Status column,
# matching the existing "done" entries already present for rows 3-4.
$ws.Range("E5").Value = "done"
$ws.Range("E6").Value = "done"
$ws.Range("E7").Value = "done"

# Move/restore the active selection to D14
$ws.Range("D14").Select()
